# Generate Report for Handback
# Updates the zh-cn (sheet 2) and de-de (sheet 3) status sheets with the
# outcome of the handback-validation run for c0f7397b-443e-4327-a2e2-4c0d285bf37b:
#   - the handback was rejected because it references a stale commit
#   - a "Latest Target File" (source .md) hyperlink + handback xlf/time are recorded
#   - an Error Detail message is recorded explaining the stale-version problem
# Also widens a few columns (I, J, P) on both sheets so the new text fits.

$wb = $excel.ActiveWorkbook

$latestMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/daff3162f76b0dc09c502e7be98ef40a804baf64/e2e/c0f7397b-443e-4327-a2e2-4c0d285bf37b.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dffd1d0dc2def5b65a146384b64cd36ac5ef6bd9/e2e/c0f7397b-443e-4327-a2e2-4c0d285bf37b.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/daff3162f76b0dc09c502e7be98ef40a804baf64/e2e/c0f7397b-443e-4327-a2e2-4c0d285bf37b.md."

function Update-HandbackRow6 {
    param($ws, [string]$handbackXlf, [string]$handbackTime)

    # Column widths: I, J and P grow to 40 characters.
    $ws.Columns.Item(9).ColumnWidth = 40
    $ws.Columns.Item(10).ColumnWidth = 40
    $ws.Columns.Item(16).ColumnWidth = 40

    # I6 -- Latest Target File: the source .md now has a target + hyperlink.
    $i6 = $ws.Cells.Item(6, 9)
    $i6.Value = "c0f7397b-443e-4327-a2e2-4c0d285bf37b.md"
    $i6.Hyperlinks.Add($i6, $latestMdUrl)
    $i6.Font.Underline = 2
    $i6.Font.Color = 15570276

    # J6 -- Latest Handback File
    $ws.Cells.Item(6, 10).Value = $handbackXlf

    # K6 -- Latest Handback DateTime
    $ws.Cells.Item(6, 11).Value = $handbackTime

    # P6 -- Error Detail
    $ws.Cells.Item(6, 16).Value = $errorDetail
}

$wsZhCn = $wb.Worksheets.Item("zh-cn")
Update-HandbackRow6 -ws $wsZhCn `
    -handbackXlf "c0f7397b-443e-4327-a2e2-4c0d285bf37b.1669576dddce01f5f43b4ecf6e6880b332255bfc.zh-cn.xlf" `
    -handbackTime "2016-10-17 14:05:02"

$wsDeDe = $wb.Worksheets.Item("de-de")
Update-HandbackRow6 -ws $wsDeDe `
    -handbackXlf "c0f7397b-443e-4327-a2e2-4c0d285bf37b.1669576dddce01f5f43b4ecf6e6880b332255bfc.de-de.xlf" `
    -handbackTime "2016-10-17 14:05:40"
